# Sprint 3 Hours Log -- Mason
# Fill in the next time-log entry (row 5) with Wednesday 2017-04-20's hours,
# and advance the active-cell selection down to row 6, exactly like a user
# would after finishing typing the row and pressing Enter.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data: date, hours, user story id, description
$ws.Range("A5").Value = 42845
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = "SF-17"
$ws.Range("D5").Value = "Created the sort method and properly implemented it"

# Re-apply the row's font so the new text picks up its own (visually
# identical) font/style entry, matching the rest of the log rows.
$ws.Range("C5:D5").Font.ThemeColor = 1

# Move the selection to A6, as if the user pressed Enter after row 5
$ws.Range("A6").Select()
